$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.040.89'
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").Value = '2.238.77'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.77'
$ws.Range("E5").Value = '  -4.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.75'
$ws.Range("E6").Value = '  -6.40%  '
$ws.Range("E7").Value = '  -1.67%  '
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("E9").Value = '  -5.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.91'
$ws.Range("E10").Value = '  -6.67%  '
$ws.Range("E11").Value = '  -3.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.23'
$ws.Range("E12").Value = '  -5.25%  '
$ws.Range("E13").Value = '  -2.75%  '
$ws.Range("D14").Value = '2.579.66'
$ws.Range("E14").Value = '  -0.87%  '
$ws.Range("D15").Value = '2.239.90'
$ws.Range("E15").Value = '  -0.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.824'
$ws.Range("E16").Value = '  -4.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.59'
$ws.Range("E17").Value = '  -5.91%  '
$ws.Range("D18").Value = '43.910.26'
$ws.Range("E18").Value = '  -0.33%  '
$ws.Range("D19").Value = '0.0₃0957'
$ws.Range("E19").Value = '  -3.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.24'
$ws.Range("E21").Value = '  -4.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '64.86'
$ws.Range("E22").Value = '  -1.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.90'
$ws.Range("E23").Value = '  +0.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.92'
$ws.Range("E24").Value = '  -7.54%  '
$ws.Range("E25").Value = '  -7.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.01'
$ws.Range("E26").Value = '  +0.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.96'
$ws.Range("E27").Value = '  -2.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.79'
$ws.Range("E28").Value = '  -1.78%  '
$ws.Range("E29").Value = '  -2.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.96'
$ws.Range("E30").Value = '  -4.32%  '
$ws.Range("E31").Value = '  -1.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.77'
$ws.Range("E32").Value = '  -4.19%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0807'
$ws.Range("E33").Value = '  -5.34%  '
$ws.Range("E34").Value = '  +8.86%  '
$ws.Range("E35").Value = '  -2.80%  '
$ws.Range("E36").Value = '  -6.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.118'
$ws.Range("E38").Value = '  -9.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.24'
$ws.Range("E39").Value = '  -8.78%  '
$ws.Range("E40").Value = '  -9.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.83'
$ws.Range("E41").Value = '  -9.12%  '
$ws.Range("E42").Value = '  -4.94%  '
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("D44").Value = '1.735.48'
$ws.Range("E44").Value = '  -2.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '85.85'
$ws.Range("E45").Value = '  +3.94%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.188'
$ws.Range("E46").Value = '  -5.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '100.09'
$ws.Range("E47").Value = '  -4.79%  '
$ws.Range("E48").Value = '  -5.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '14.66'
$ws.Range("E49").Value = '  +1.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '69.00'
$ws.Range("E50").Value = '  -8.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.07'
$ws.Range("E51").Value = '  -4.28%  '
